$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# MVPFeatures sheet: insert 6 new rows at the top (below the header row) for
# the new "ESG reporting pivot" thinking / MVP notes, then activate this
# sheet (it becomes the tab the workbook opens to).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("MVPFeatures")
$dailyNotes = $wb.Worksheets.Item("DailyNotes")

# Insert 6 blank rows above the old row 2 - formatting is inherited from row 1
# (A:s5 bold/underline/wrap, B:s9 bold), matching the final file exactly for
# rows 3 (B col), 4-6 (B col) and 7 (A+B col), so only a handful of cells need
# explicit value/format work below.
$ws.Rows("2:7").Insert()

# Row 2: new dated heading entry (style copied from the existing "Other
# possible activities in future version:" heading cell and the existing
# bold-dated "Assignment / Deadline" style cell so we reuse the same cellXfs
# entries rather than inventing new ones).
$ws.Range("A2").Value = "'-- Given need for ESG reporting in 2024, there will be fast, massive adoption of need to know. Pivot away from just renewable energy research. Start with collecting/summarizing data and step 1a of Decarb Strategist Process"
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Value = "10/1/2023"
$dailyNotes.Range("A51").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null

# Row 3: bold + underlined sub-heading.
$ws.Range("A3").Value = "'-- 3 PERTINENT features for MVP:"
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Underline = $true
$ws.Range("A3").WrapText = $true

# Rows 4-6: the three MVP feature bullets (plain, wrapped text).
$ws.Range("A4").Value = "'--- (1) Collect reported ESG data from any company that publishes these metrics and summarizes what's used"
$ws.Range("A4").WrapText = $true

$ws.Range("A5").Value = "'--- (2) Sends users related, real-time updates from the web based on past queries"
$ws.Range("A5").WrapText = $true

$ws.Range("A6").Value = "'--- (3) Walks users through process of collecting internal data for comparison"
$ws.Range("A6").WrapText = $true

# Selection lands on the new blank spacer row (A7), matching the saved file,
# and this sheet becomes the active / tab-selected sheet in the workbook.
$ws.Range("A7").Select()
$ws.Activate()

$wb.Application.ActiveWindow.Zoom = $wb.Application.ActiveWindow.Zoom
